# Add an "hour" column to the "combo" sheet, inserted between the
# existing "time" column (B) and the "temp_mean" column (old C, now D).
#
# This mirrors a manual edit in Excel: select column C, Insert a new
# column (pushing temp_mean/temp_se/humid_mean/humid_se one slot to the
# right), fill in a header, then fill down an elapsed-hours formula for
# each of the two logged groups (tank_A readings every 15 min, tank_B
# readings every 10 min with a longer overnight gap between rows 31/32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("combo")

# --- insert the new column -------------------------------------------------
$ws.Columns.Item(3).Insert()

# Header
$ws.Range("C1").Value = "hour"

# --- first block (tank_A group, rows 2-29): 15-minute sampling interval ---
$ws.Range("C2").Value = 0
$ws.Range("C3").Formula = "=C2+(1*15/60)"
$ws.Range("C4").Formula = "=C3+(1*15/60)"
$ws.Range("C5").Formula = "=C4+(1*15/60)"
$ws.Range("C6").Formula = "=C5+(1*15/60)"
$ws.Range("C7:C29").Formula = "=C6+(1*15/60)"

# --- second block (tank_B group, rows 30-65): 10-minute sampling interval,
#     with an overnight (+1 hour-on-clock... +1 day) gap between the first
#     two readings and the rest of the run ---
$ws.Range("C30").Value = 0
$ws.Range("C31").Formula = "=C30+(1*10/60)"
$ws.Range("C32").Formula = "=C31+1+(1*10/60)"
$ws.Range("C33:C65").Formula = "=C32+(1*10/60)"

# --- cosmetic: give the new column roughly the same width/format as the
#     neighbouring "time" column it was inserted next to ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- restore the author's cursor position on this sheet ---
$ws.Activate()
$ws.Range("C29").Select()
